# Append the new 2025-04-01 price row (row 31) to every price sheet,
# carrying forward the last known (2025-03-31) price, mirroring how the
# prior rows in each sheet were populated (literal text in column A/B).

$wb = $excel.ActiveWorkbook

$sheetValues = @{
    "N-Dense"                   = "40"
    "N-Type"                    = "42"
    "N-type Wafer"               = "1.25"
    "Cell Topcon 183mm"         = "0.303"
    "Module Topcon 183mm"       = "0.1"
    "Silver Rear_side"          = "5,509"
    "Silver Busbar front-side"  = "8,247"
    "Silver finger front-side"  = "8,297"
    "USD_CNY"                   = "7.2752"
}

$newDate = "2025-04-01"

foreach ($name in $sheetValues.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $newRow = $ws.Cells.Item(30, 1).Row + 1

    $dateCell = $ws.Cells.Item($newRow, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newDate

    $priceCell = $ws.Cells.Item($newRow, 2)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $sheetValues[$name]
}
